# Insert a new daily-stats row right before the existing "2026/12/29" block
# (old row 679), shifting every subsequent row down by one. The new row
# carries the date 2026/01/22 (day-of-week 木) with time-slot 23 and
# ranking 15 -- matching an extra sample for a date that already appears
# earlier in the sheet (row 678).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 679..720 down to 680..721, leaving a blank row at 679.
$ws.Rows.Item(679).Insert()

# Fill the new row 679. Force text format on the date cell first so Excel
# doesn't auto-coerce the "yyyy/mm/dd" string into a real date serial
# (the rest of the column is stored as plain text), then clear the
# formatting override so the cell ends up with the same "no explicit
# style" look as its neighbours.
$ws.Range("A679").NumberFormat = "@"
$ws.Range("A679").Value = "2026/01/22"
$ws.Range("A679").ClearFormats()

$ws.Range("B679").Value = "木"
$ws.Range("C679").Value = 23
$ws.Range("D679").Value = 15
